# Updates division problems in the single table of the document.
$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

$cell = $table.Cell(1, 1)
$cellRange = $cell.Range
$cellRange.MoveEnd(1, -1) | Out-Null  # exclude end-of-cell marker
if ($cellRange.Text -ne '76÷9=8, 4') { throw "Unexpected cell text at (1,1): $($cellRange.Text)" }
$cellRange.Text = '37÷4=9, 1'

$cell = $table.Cell(1, 2)
$cellRange = $cell.Range
$cellRange.MoveEnd(1, -1) | Out-Null  # exclude end-of-cell marker
if ($cellRange.Text -ne '48÷8=6, 0') { throw "Unexpected cell text at (1,2): $($cellRange.Text)" }
$cellRange.Text = '73÷9=8, 1'

$cell = $table.Cell(1, 3)
$cellRange = $cell.Range
$cellRange.MoveEnd(1, -1) | Out-Null  # exclude end-of-cell marker
if ($cellRange.Text -ne '46÷8=5, 6') { throw "Unexpected cell text at (1,3): $($cellRange.Text)" }
$cellRange.Text = '14÷9=1, 5'

$cell = $table.Cell(1, 4)
$cellRange = $cell.Range
$cellRange.MoveEnd(1, -1) | Out-Null  # exclude end-of-cell marker
if ($cellRange.Text -ne '98÷8=12, 2') { throw "Unexpected cell text at (1,4): $($cellRange.Text)" }
$cellRange.Text = '43÷4=10, 3'

$cell = $table.Cell(1, 5)
$cellRange = $cell.Range
$cellRange.MoveEnd(1, -1) | Out-Null  # exclude end-of-cell marker
if ($cellRange.Text -ne '56÷6=9, 2') { throw "Unexpected cell text at (1,5): $($cellRange.Text)" }
$cellRange.Text = '66÷8=8, 2'

$cell = $table.Cell(5, 1)
$cellRange = $cell.Range
$cellRange.MoveEnd(1, -1) | Out-Null  # exclude end-of-cell marker
if ($cellRange.Text -ne '62÷2=31, 0') { throw "Unexpected cell text at (5,1): $($cellRange.Text)" }
$cellRange.Text = '19÷2=9, 1'

$cell = $table.Cell(5, 2)
$cellRange = $cell.Range
$cellRange.MoveEnd(1, -1) | Out-Null  # exclude end-of-cell marker
if ($cellRange.Text -ne '74÷3=24, 2') { throw "Unexpected cell text at (5,2): $($cellRange.Text)" }
$cellRange.Text = '29÷4=7, 1'

$cell = $table.Cell(5, 3)
$cellRange = $cell.Range
$cellRange.MoveEnd(1, -1) | Out-Null  # exclude end-of-cell marker
if ($cellRange.Text -ne '65÷2=32, 1') { throw "Unexpected cell text at (5,3): $($cellRange.Text)" }
$cellRange.Text = '49÷5=9, 4'

$cell = $table.Cell(5, 4)
$cellRange = $cell.Range
$cellRange.MoveEnd(1, -1) | Out-Null  # exclude end-of-cell marker
if ($cellRange.Text -ne '62÷7=8, 6') { throw "Unexpected cell text at (5,4): $($cellRange.Text)" }
$cellRange.Text = '87÷9=9, 6'

$cell = $table.Cell(5, 5)
$cellRange = $cell.Range
$cellRange.MoveEnd(1, -1) | Out-Null  # exclude end-of-cell marker
if ($cellRange.Text -ne '42÷2=21, 0') { throw "Unexpected cell text at (5,5): $($cellRange.Text)" }
$cellRange.Text = '34÷2=17, 0'

$cell = $table.Cell(9, 1)
$cellRange = $cell.Range
$cellRange.MoveEnd(1, -1) | Out-Null  # exclude end-of-cell marker
if ($cellRange.Text -ne '53÷3=17, 2') { throw "Unexpected cell text at (9,1): $($cellRange.Text)" }
$cellRange.Text = '97÷9=10, 7'

$cell = $table.Cell(9, 2)
$cellRange = $cell.Range
$cellRange.MoveEnd(1, -1) | Out-Null  # exclude end-of-cell marker
if ($cellRange.Text -ne '87÷3=29, 0') { throw "Unexpected cell text at (9,2): $($cellRange.Text)" }
$cellRange.Text = '35÷8=4, 3'

$cell = $table.Cell(9, 3)
$cellRange = $cell.Range
$cellRange.MoveEnd(1, -1) | Out-Null  # exclude end-of-cell marker
if ($cellRange.Text -ne '18÷8=2, 2') { throw "Unexpected cell text at (9,3): $($cellRange.Text)" }
$cellRange.Text = '52÷5=10, 2'

$cell = $table.Cell(9, 4)
$cellRange = $cell.Range
$cellRange.MoveEnd(1, -1) | Out-Null  # exclude end-of-cell marker
if ($cellRange.Text -ne '90÷8=11, 2') { throw "Unexpected cell text at (9,4): $($cellRange.Text)" }
$cellRange.Text = '62÷9=6, 8'

$cell = $table.Cell(9, 5)
$cellRange = $cell.Range
$cellRange.MoveEnd(1, -1) | Out-Null  # exclude end-of-cell marker
if ($cellRange.Text -ne '13÷9=1, 4') { throw "Unexpected cell text at (9,5): $($cellRange.Text)" }
$cellRange.Text = '41÷8=5, 1'

$cell = $table.Cell(13, 1)
$cellRange = $cell.Range
$cellRange.MoveEnd(1, -1) | Out-Null  # exclude end-of-cell marker
if ($cellRange.Text -ne '32÷3=10, 2') { throw "Unexpected cell text at (13,1): $($cellRange.Text)" }
$cellRange.Text = '84÷5=16, 4'

$cell = $table.Cell(13, 2)
$cellRange = $cell.Range
$cellRange.MoveEnd(1, -1) | Out-Null  # exclude end-of-cell marker
if ($cellRange.Text -ne '25÷6=4, 1') { throw "Unexpected cell text at (13,2): $($cellRange.Text)" }
$cellRange.Text = '17÷3=5, 2'

$cell = $table.Cell(13, 3)
$cellRange = $cell.Range
$cellRange.MoveEnd(1, -1) | Out-Null  # exclude end-of-cell marker
if ($cellRange.Text -ne '78÷4=19, 2') { throw "Unexpected cell text at (13,3): $($cellRange.Text)" }
$cellRange.Text = '29÷6=4, 5'

$cell = $table.Cell(13, 4)
$cellRange = $cell.Range
$cellRange.MoveEnd(1, -1) | Out-Null  # exclude end-of-cell marker
if ($cellRange.Text -ne '77÷7=11, 0') { throw "Unexpected cell text at (13,4): $($cellRange.Text)" }
$cellRange.Text = '84÷8=10, 4'

$cell = $table.Cell(13, 5)
$cellRange = $cell.Range
$cellRange.MoveEnd(1, -1) | Out-Null  # exclude end-of-cell marker
if ($cellRange.Text -ne '72÷3=24, 0') { throw "Unexpected cell text at (13,5): $($cellRange.Text)" }
$cellRange.Text = '34÷4=8, 2'

$cell = $table.Cell(17, 1)
$cellRange = $cell.Range
$cellRange.MoveEnd(1, -1) | Out-Null  # exclude end-of-cell marker
if ($cellRange.Text -ne '85÷5=17, 0') { throw "Unexpected cell text at (17,1): $($cellRange.Text)" }
$cellRange.Text = '62÷9=6, 8'

$cell = $table.Cell(17, 2)
$cellRange = $cell.Range
$cellRange.MoveEnd(1, -1) | Out-Null  # exclude end-of-cell marker
if ($cellRange.Text -ne '64÷4=16, 0') { throw "Unexpected cell text at (17,2): $($cellRange.Text)" }
$cellRange.Text = '79÷3=26, 1'

$cell = $table.Cell(17, 3)
$cellRange = $cell.Range
$cellRange.MoveEnd(1, -1) | Out-Null  # exclude end-of-cell marker
if ($cellRange.Text -ne '35÷7=5, 0') { throw "Unexpected cell text at (17,3): $($cellRange.Text)" }
$cellRange.Text = '37÷2=18, 1'

$cell = $table.Cell(17, 4)
$cellRange = $cell.Range
$cellRange.MoveEnd(1, -1) | Out-Null  # exclude end-of-cell marker
if ($cellRange.Text -ne '50÷2=25, 0') { throw "Unexpected cell text at (17,4): $($cellRange.Text)" }
$cellRange.Text = '43÷9=4, 7'

$cell = $table.Cell(17, 5)
$cellRange = $cell.Range
$cellRange.MoveEnd(1, -1) | Out-Null  # exclude end-of-cell marker
if ($cellRange.Text -ne '78÷4=19, 2') { throw "Unexpected cell text at (17,5): $($cellRange.Text)" }
$cellRange.Text = '51÷9=5, 6'
